$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header block - cardholder name and card number
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 01.06.2025"

# Transaction rows
$ws.Range("B6").Value = "03.06."
$ws.Range("C6").Value = "04.06."
$ws.Range("D6").Value = "BEITRAG Allianz SE K-45725151"
$ws.Range("E6").Value = "54,62-"

$ws.Range("B7").Value = "06.06."
$ws.Range("C7").Value = "07.06."
$ws.Range("D7").Value = "EBAY MKTPLC EU TRDJOD"
$ws.Range("E7").Value = "93,05-"

$ws.Range("B8").Value = "08.06."
$ws.Range("C8").Value = "09.06."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 46721599"
$ws.Range("E8").Value = "38,87-"

$ws.Range("B9").Value = "11.06."
$ws.Range("C9").Value = "12.06."
$ws.Range("D9").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E9").Value = "24,64-"

$ws.Range("B10").Value = "13.06."
$ws.Range("C10").Value = "14.06."
$ws.Range("D10").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 98041088"
$ws.Range("E10").Value = "87,31-"

$ws.Range("B11").Value = "17.06."
$ws.Range("C11").Value = "18.06."
$ws.Range("D11").Value = "MCDONALDS Saarlouis"
$ws.Range("E11").Value = "25,61-"

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 22.06.2025"
$ws.Range("E12").Value = "324,10-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 02.07.2025"
